$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Password" value in G2 (shared string content change)
$ws.Range("G2").Value = "aH5o@UmNP5"

# Widen the new column G so its generated text fits (closest attainable
# value to the authored 17.81640625 given this host's 1/6-character grid)
$ws.Range("G:G").ColumnWidth = 17

# Move/restore the active selection to G2
$ws.Range("G2").Select()
